$wb = $excel.ActiveWorkbook

# Workbook window size change
$excel.Width = 8860
$excel.Height = 3290

# Sheet1: scroll + selection change, and clear row 9 contents
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Activate()
$ws1.Application.ActiveWindow.ScrollRow = 7
$ws1.Range("A9:E9").Select()
$ws1.Range("A9:E9").ClearContents()

# Sheet2: remove the now-redundant separate column width for column A
$ws2 = $wb.Worksheets.Item("ESRI_MAPINFO_SHEET")
$ws2.Columns.Item(1).ColumnWidth = $ws2.Columns.Item(2).ColumnWidth
